# Fix map farm locations so that everything is in or close to King county.
#
# A handful of farms had coordinates far outside King County (typos /
# mis-geocoded addresses). Replace those coordinate values with corrected
# ones that are in (or close to) King County. The worksheet is sorted
# alphabetically by farm name in column A, with coordinates in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of Farm Name -> corrected "lat,long" coordinate string.
$fixes = @{
    "Caruso Farm"         = "47.8782821,-122.0526352"
    "Flying Tomato Farm"  = "47.5555656,-122.3209544"
    "High & Dry Farm"     = "47.8877264,-121.7938542"
    "Hima Farms"          = "47.6757423,-122.2394706"
    "Local Roots Farm"    = "47.6349728,-122.0819622"
    "Lowlands Farm"       = "47.8062607,-122.122593"
    "NW Farms"            = "47.3239497,-122.2772907"
    "Orange Star Farm"    = "47.8634311,-121.9435757"
    "Oxbow Farm"          = "47.6777367,-122.0058828"
    "Radicle Roots Farm"  = "47.7605015,-122.3874161"
    "Skylight Farms"      = "47.836134,-122.0852195"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ($fixes.ContainsKey($name)) {
        $ws.Cells.Item($r, 2).Value = $fixes[$name]
    }
}

# Update the saved view: scroll the window so row 6 is the top-left visible
# row, and select B16 as the active cell (matches the author's last editing
# position before committing).
$ws.Range("B16").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
